$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Desktop Computer" (DKS) rows 8-10, which shifts the rows below
# (CMR, SCN, PRT) up. Excel also rebuilds the shared-strings table on save,
# dropping the now-unused DKS/Dekstop/Desktop Computer/... entries.
$ws.Range("A8:G10").EntireRow.Delete()

# Apply the print/page setup that was left on the sheet
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

# Leave the active selection on E10, matching the final state left in the file
$ws.Range("E10").Select() | Out-Null
